$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) - update "想去人数" (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1601
$wsExpo.Range("F4").Value = 8546
$wsExpo.Range("F6").Value = 65
$wsExpo.Range("F7").Value = 24
$wsExpo.Range("F10").Value = 108
$wsExpo.Range("F13").Value = 9201
$wsExpo.Range("F16").Value = 214
$wsExpo.Range("F18").Value = 343
$wsExpo.Range("F19").Value = 6092
$wsExpo.Range("F20").Value = 1045
$wsExpo.Range("F21").Value = 62

# Sheet "全部类型" (all types) - same events, shifted row numbers, update matching values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1601
$wsAll.Range("F4").Value = 8546
$wsAll.Range("F6").Value = 65
$wsAll.Range("F7").Value = 24
$wsAll.Range("F10").Value = 108
$wsAll.Range("F15").Value = 9201
$wsAll.Range("F18").Value = 214
$wsAll.Range("F20").Value = 343
$wsAll.Range("F21").Value = 6092
$wsAll.Range("F22").Value = 1045
$wsAll.Range("F23").Value = 62
